$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.113764762878418
$ws.Range("B1").Value = 2.446182727813721
$ws.Range("C1").Value = 5.22813606262207
$ws.Range("D1").Value = 2.205544948577881
$ws.Range("E1").Value = 1.271423697471619
